$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 489, shifting rows 489:540 down to 490:541
$ws.Rows.Item(489).Insert()

# Populate the new row 489 with the new record.
# Columns A,B,C,E,F,G,H,I,N,Q,R are unchanged (same as the rest of the block),
# while D,J,K,L,M,O,P carry new values.
$ws.Cells.Item(489, 1).Value = 10
$ws.Cells.Item(489, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(489, 3).Value = "La Araucanía"
$ws.Cells.Item(489, 4).Value = 44946
$ws.Cells.Item(489, 5).Value = 9
$ws.Cells.Item(489, 6).Value = 100112037
$ws.Cells.Item(489, 7).Value = "Cebollín"
$ws.Cells.Item(489, 8).Value = "Sin especificar"
$ws.Cells.Item(489, 9).Value = "Primera"
$ws.Cells.Item(489, 10).Value = 125
$ws.Cells.Item(489, 11).Value = 6000
$ws.Cells.Item(489, 12).Value = 6000
$ws.Cells.Item(489, 13).Value = 6000
$ws.Cells.Item(489, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(489, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(489, 16).Value = 500
$ws.Cells.Item(489, 17).Value = 12
$ws.Cells.Item(489, 18).Value = "Hortaliza"
